$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 data values (environment block) ---
# A2 keeps its quote-prefixed style (s=1); restore it explicitly via a format-only
# paste, since assigning .Value/.Style directly re-evaluates (and drops) quotePrefix.
$ws.Range("A2").Value = "i-preproducciongestion.segurossura.com.ar"
$ws.Range("Q2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D2").Value = "silverarrow"

# Update the hyperlink cell text + underlying hyperlink target
$b2Style = $ws.Range("B2").Style
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do")
$ws.Range("B2").Value = "https://i-preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("B2").Style = $b2Style

# Keep policy number a text value (avoid losing the leading zero)
$ws.Range("E2").Value = "'04104013566"

# --- New column F: "Anulada" / "N" ---
$ws.Range("F1").Value = "Anulada"
$ws.Range("F2").Value = "N"

# --- Formatting-only tweaks on empty helper cells in column E ---
$ws.Range("E4").Copy()
$ws.Range("E3").PasteSpecial(-4122)

$ws.Range("E8").Copy()
$ws.Range("E7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Remove row 11 (rows 12-16 shift up to 11-15) ---
$ws.Rows("11:11").Delete()

# --- Final selection as recorded in the workbook ---
$ws.Range("F3").Select()
